$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "RS Activist" task row (A6) needs a small text fix.
$ws.Range("A6").Value = "RS Activist -"
